$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new rows at 13-14 for the professor names (shifts everything below down by 2)
$ws.Rows("13:14").Insert()

# Remove the stray empty A13/A14 cells created by the row insert's format inheritance
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# 2) Fill in the new professor rows (column B and C only, no column A label)
$ws.Range("B13").Value = '5840963 - Daniela Camargo Vernilli'
$ws.Range("C13").Value = '5840963 - Daniela Camargo Vernilli'
$ws.Range("B14").Value = '5840820 - Gustavo Aristides Santana Martinez'
$ws.Range("C14").Value = '5840820 - Gustavo Aristides Santana Martinez'

# Re-apply the correct column styles (B=style 2, C=style 3) lost when the row was cleared
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Update the text content that changed
$ws.Range("B10").Value = 'Munir o aluno de ferramentas para: especificar materiais para aplicação na indústria química e definir os meios possíveis de processar os materiais comuns a sua área de atuação.'
$ws.Range("C10").Value = 'Munir o aluno de ferramentas para: especificar materiais para aplicação na indústria química e definir os meios possíveis de processar os materiais comuns a sua área de atuação.'

$ws.Range("B15").Value = 'Introdução aos materiais para a indústria química e outras, propriedades,  especificações, seleção, fabricação, aplicação e possíveis falhas.'
$ws.Range("C15").Value = 'Introdução aos materiais para a indústria química e outras, propriedades,  especificações, seleção, fabricação, aplicação e possíveis falhas.'

$ws.Range("B17").Value = 'Introdução aos materiais. - Seleção de materiais. - Fatores que influenciam na seleção dos materiais (indústria química, petroquímica, Nuclear e outras), melhoria das propriedades mecânicas dos metais. - Falhas em serviço e em processo.  Produtos siderúrgicos para aplicação em indústrias químicas - Aços carbono e especiais - Ferro fundido. – Processo de fabricação de aços e ferros fundidos, especificações, propriedades e aplicações.  Metais e ligas não ferrosas: especificações, propriedades e aplicações. Introdução à corrosão. - Causas e formas de corrosão. Proteção de superfícies metálicas contra a corrosão, revestimentos. Requisitos específicos de materiais metálicos para a indústria de óleo e gás.  Materiais não metálicos. Especificações, propriedades e aplicações.'
$ws.Range("C17").Value = 'Introdução aos materiais. - Seleção de materiais. - Fatores que influenciam na seleção dos materiais (indústria química, petroquímica, Nuclear e outras), melhoria das propriedades mecânicas dos metais. - Falhas em serviço e em processo.  Produtos siderúrgicos para aplicação em indústrias químicas - Aços carbono e especiais - Ferro fundido. – Processo de fabricação de aços e ferros fundidos, especificações, propriedades e aplicações.  Metais e ligas não ferrosas: especificações, propriedades e aplicações. Introdução à corrosão. - Causas e formas de corrosão. Proteção de superfícies metálicas contra a corrosão, revestimentos. Requisitos específicos de materiais metálicos para a indústria de óleo e gás.  Materiais não metálicos. Especificações, propriedades e aplicações.'

$ws.Range("B20").Value = 'Duas provas'
$ws.Range("C20").Value = 'Duas provas'

$ws.Range("B21").Value = 'Serão aplicadas duas avaliações (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF= (P1+P2)/2'
$ws.Range("C21").Value = 'Serão aplicadas duas avaliações (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF= (P1+P2)/2'

$ws.Range("B22").Value = 'Para o aluno que obtiver Nota Final maior ou igual a 3,0 e menor do que 5,0, será aplicada uma avaliação de recuperação (NR), com pontuação de 0 a 10, que levará ao cálculo da média final(MF) através da seguinte expressão: MF=(NF+NR)/2.onde: NF=Nota Final e NR=Nota da Prova de Recuperação.'
$ws.Range("C22").Value = 'Para o aluno que obtiver Nota Final maior ou igual a 3,0 e menor do que 5,0, será aplicada uma avaliação de recuperação (NR), com pontuação de 0 a 10, que levará ao cálculo da média final(MF) através da seguinte expressão: MF=(NF+NR)/2.onde: NF=Nota Final e NR=Nota da Prova de Recuperação.'

$ws.Range("B23").Value = '1)Telles, P. C. S. - Materiais para Equipamentos de Processo - Ed. Interciência, 4º Ed., 1989.2)Bresciani, F., E. - Seleção de Materiais Metálicos - Ed. da UNICAMP, 2º Ed.3)Freire, J. M. -Materiais de Construção Mecânica - Ed. Livros Técnicos e Científicos, Editora 1993.4)A. Remy/ M. Gay/ R. Gonthier - Materiais - Hemus Editora Limitada - 2ª Edição.5)Chiaverini, V.Tecnologia Mecânica - Materiais de Construção Mecânica - Vol. II - Ed. McGraw Hill do Brasil Ltda.6)Gentil, V. - Corrosão. - Ed. Guanabara Dois, 1982.'
$ws.Range("C23").Value = '1)Telles, P. C. S. - Materiais para Equipamentos de Processo - Ed. Interciência, 4º Ed., 1989.2)Bresciani, F., E. - Seleção de Materiais Metálicos - Ed. da UNICAMP, 2º Ed.3)Freire, J. M. -Materiais de Construção Mecânica - Ed. Livros Técnicos e Científicos, Editora 1993.4)A. Remy/ M. Gay/ R. Gonthier - Materiais - Hemus Editora Limitada - 2ª Edição.5)Chiaverini, V.Tecnologia Mecânica - Materiais de Construção Mecânica - Vol. II - Ed. McGraw Hill do Brasil Ltda.6)Gentil, V. - Corrosão. - Ed. Guanabara Dois, 1982.'

# 4) Fix the column definitions: column A should only be width 30.7109375 on its own (max=1), not overlapping column B
$ws.Columns("A").ColumnWidth = $ws.Columns("A").ColumnWidth
